$wb = $excel.ActiveWorkbook

# Sheets 1-5 all have a header row (row 1) where E1 was mistakenly left as a
# stray numeric placeholder (688.2298852957874) instead of the intended
# year/period label. Fix the label on each sheet. Sheets 1,2,3,5 use plain
# "2050"; formatting the cell as Text first keeps Excel's smart-typing from
# re-interpreting the digits as a number. Sheet 4 uses year-range labels
# elsewhere (e.g. "2015-2030"), so its fixed label is "2041-2050" - already
# non-numeric text, no special formatting needed.
#
# Sheets 1-4 (and 6) also have a stray "Total" row at the bottom of their
# tables that needs to be removed entirely.

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E1").NumberFormat = "@"
$ws1.Range("E1").Value = "2050"
$ws1.Rows.Item(13).Delete()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("E1").NumberFormat = "@"
$ws2.Range("E1").Value = "2050"
$ws2.Rows.Item(13).Delete()

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("E1").NumberFormat = "@"
$ws3.Range("E1").Value = "2050"
$ws3.Rows.Item(13).Delete()

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("E1").Value = "2041-2050"
$ws4.Rows.Item(13).Delete()

$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("E1").NumberFormat = "@"
$ws5.Range("E1").Value = "2050"

$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
